$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for shared-string table placement: E1, then E3, then E2
$ws.Range("E1").Value = "23-123146"
$ws.Range("E3").Value = "23-123148"
$ws.Range("E2").Value = "23-123137"

$ws.Range("D3").Value2 = 44999

$ws.Range("E2").Select()
